$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCDTtiNTY")

# C7: replace literal 0 with formula =B2
$ws.Range("C7").Formula = "=B2"

# D7: add new note text
$ws.Range("D7").Value = "motobikes F are assumed same with LDV P"

# Make the SoCDTtiNTY sheet the active (selected) sheet/tab
$ws.Activate()
$ws.Select()

# Set the selection on SoCDTtiNTY sheet to E15
$ws.Range("E15").Select()
